# "run modelv5 at all depths"
#
# The model output now covers every depth, the row labels moved to the
# longer "modelv5..." names, and the sheet was reviewed with the focus
# down near row 53 instead of the old single-cell selection at J62.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A so the longer modelv5 row labels aren't truncated.
$ws.Columns.Item(1).ColumnWidth = 16

# Scroll the window back up so row 34 is the first visible row.
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1

# Select row 53 in its entirety - the row of interest for the new run.
[void]$ws.Rows.Item(53).Select()

# The workbook window was maximized while reviewing the new results.
$excel.ActiveWindow.WindowState = -4137
